# Fixed some errors in Aeordynamic and Stability
# Updates recomputed DOC (Direct Operating Cost) figures on the single
# "DOC" worksheet. Only numeric value cells change; labels/formatting
# are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Depreciation
$ws.Range("B2").Value = 920.8819664592633
$ws.Range("D2").Value = 4.604409832296318
$ws.Range("E2").Value = 6.395013655967108

# Row 3 - Interest
$ws.Range("B3").Value = 690.6614748444475
$ws.Range("D3").Value = 3.4533073742222387
$ws.Range("E3").Value = 4.7962602419753315

# Row 4 - Insurance
$ws.Range("B4").Value = 122.78426219456844
$ws.Range("D4").Value = 0.6139213109728424
$ws.Range("E4").Value = 0.8526684874622811

# Row 5 - DOC Capital
$ws.Range("B5").Value = 1734.3277034982789
$ws.Range("D5").Value = 8.671638517491397
$ws.Range("E5").Value = 12.043942385404717

# Row 7 - Cockpit Crew
$ws.Range("B7").Value = 254.16342274275667
$ws.Range("D7").Value = 1.2708171137137836
$ws.Range("E7").Value = 1.765023769046922

# Row 8 - Cabin Crew
$ws.Range("B8").Value = 84.72114091425223
$ws.Range("C8").Value = 40.0
$ws.Range("D8").Value = 0.4236057045712613
$ws.Range("E8").Value = 0.5883412563489739

# Row 9 - DOC Crew
$ws.Range("B9").Value = 338.8845636570089
$ws.Range("C9").Value = 160.0
$ws.Range("D9").Value = 1.6944228182850452
$ws.Range("E9").Value = 2.3533650253958958

# Row 11 - DOC Fuel
$ws.Range("B11").Value = 252.828738682599
$ws.Range("C11").Value = 119.36984604043116
$ws.Range("D11").Value = 1.2641436934129955
$ws.Range("E11").Value = 1.7557551297402716

# Row 13 - Navigation charges
$ws.Range("B13").Value = 172.87210948661146
$ws.Range("C13").Value = 81.6193491358094
$ws.Range("D13").Value = 0.8643605474330576
$ws.Range("E13").Value = 1.200500760323691

# Row 18 - DOC Charges
$ws.Range("B18").Value = 109.77378952399826
$ws.Range("C18").Value = 51.82828670123896
$ws.Range("D18").Value = 0.5488689476199915
$ws.Range("E18").Value = 0.7623179828055437

# Row 20 - Airframe Maintenance Charges
$ws.Range("B20").Value = 394.70554059765806
$ws.Range("C20").Value = 186.3551582701874
$ws.Range("D20").Value = 1.973527702988291
$ws.Range("E20").Value = 2.7410106985948484

# Row 21 - Engine Maintenance Charges
$ws.Range("B21").Value = 279.0567432187649
$ws.Range("C21").Value = 131.75306196653003
$ws.Range("D21").Value = 1.3952837160938247
$ws.Range("E21").Value = 1.9378940501303121

# Row 22 - DOC Maintenance
$ws.Range("B22").Value = 758.3539138088424
$ws.Range("C22").Value = 358.04707331615657
$ws.Range("D22").Value = 3.791769569044213
$ws.Range("E22").Value = 5.2663466236725185

# Row 25 - Total DOC
$ws.Range("B25").Value = 3194.168709170728
$ws.Range("C25").Value = 1508.0857857679716
$ws.Range("D25").Value = 15.970843545853645
$ws.Range("E25").Value = 22.18172714701895

# Row 27 - Cash DOC
$ws.Range("B27").Value = 1459.8410056724488
$ws.Range("C27").Value = 689.2452060578267
$ws.Range("D27").Value = 7.299205028362246
$ws.Range("E27").Value = 10.137784761614231
